$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44232
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 3000
$ws.Range("O2").Value = 3000
$ws.Range("P2").Value = 3000
$ws.Range("S2").Value = 1500

# Row 3
$ws.Range("D3").Value = 44194
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 3000
$ws.Range("O3").Value = 3000
$ws.Range("P3").Value = 3000
$ws.Range("R3").Value = "Provincia de Linares"
$ws.Range("S3").Value = 1500

# Row 4
$ws.Range("D4").Value = 44238
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 3600
$ws.Range("O4").Value = 4000
$ws.Range("P4").Value = 3800
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 1900

# Row 5
$ws.Range("D5").Value = 44208
$ws.Range("M5").Value = 85
$ws.Range("N5").Value = 3000
$ws.Range("O5").Value = 3000
$ws.Range("P5").Value = 3000
$ws.Range("R5").Value = "Provincia de Linares"
$ws.Range("S5").Value = 1500

# Row 6
$ws.Range("D6").Value = 44617
$ws.Range("M6").Value = 90
$ws.Range("N6").Value = 6500
$ws.Range("O6").Value = 6500
$ws.Range("P6").Value = 6500
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value = 3250

# Row 7
$ws.Range("D7").Value = 44174
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 3200
$ws.Range("O7").Value = 3200
$ws.Range("P7").Value = 3200
$ws.Range("S7").Value = 1600

# Row 8
$ws.Range("D8").Value = 44168
$ws.Range("M8").Value = 170
$ws.Range("N8").Value = 8000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 8000
$ws.Range("S8").Value = 4000

# Row 9
$ws.Range("D9").Value = 44188
$ws.Range("M9").Value = 150
$ws.Range("N9").Value = 3000
$ws.Range("O9").Value = 3400
$ws.Range("P9").Value = 3240
$ws.Range("R9").Value = "Provincia de Linares"
$ws.Range("S9").Value = 1620

# Row 10
$ws.Range("D10").Value = 44533
$ws.Range("M10").Value = 150
$ws.Range("N10").Value = 4000
$ws.Range("O10").Value = 4000
$ws.Range("P10").Value = 4000
$ws.Range("S10").Value = 2000

# Row 11
$ws.Range("D11").Value = 44586
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 5000
$ws.Range("O11").Value = 5000
$ws.Range("P11").Value = 5000
$ws.Range("S11").Value = 2500

# Row 12
$ws.Range("D12").Value = 44582
$ws.Range("M12").Value = 380
$ws.Range("N12").Value = 5000
$ws.Range("O12").Value = 5000
$ws.Range("P12").Value = 5000
$ws.Range("R12").Value = "Provincia de Curicó"
$ws.Range("S12").Value = 2500

# Row 13
$ws.Range("D13").Value = 44236
$ws.Range("M13").Value = 300
$ws.Range("N13").Value = 3600
$ws.Range("O13").Value = 4000
$ws.Range("P13").Value = 3800
$ws.Range("S13").Value = 1900

# Row 14
$ws.Range("D14").Value = 44231
$ws.Range("M14").Value = 150
$ws.Range("N14").Value = 3400
$ws.Range("O14").Value = 3400
$ws.Range("P14").Value = 3400
$ws.Range("S14").Value = 1700

# Row 15
$ws.Range("D15").Value = 44237
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 3600
$ws.Range("O15").Value = 4000
$ws.Range("P15").Value = 3800
$ws.Range("S15").Value = 1900
